# Insert a new record row at row 16 (weekly update for "Hortaliza, Terminal
# La Palmera de La Serena - Haba"), pushing the existing rows 16-44 down to
# 17-45 and growing the used range from A1:R44 to A1:R45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16..44 down to 17..45 (copies formatting/styles along with it).
$ws.Rows.Item(16).Insert()

# The new row 16 is another "Haba" record for the same market/region, so
# copy the constant columns from the row right below it (the old row 16,
# now at row 17) and set the columns that actually carry new data.
$ws.Range("A16").Value = $ws.Range("A17").Value2
$ws.Range("B16").Value = $ws.Range("B17").Value2
$ws.Range("C16").Value = $ws.Range("C17").Value2
$ws.Range("D16").Value = 45044
$ws.Range("E16").Value = $ws.Range("E17").Value2
$ws.Range("F16").Value = $ws.Range("F17").Value2
$ws.Range("G16").Value = $ws.Range("G17").Value2
$ws.Range("H16").Value = $ws.Range("H17").Value2
$ws.Range("I16").Value = $ws.Range("I17").Value2
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("N16").Value = $ws.Range("N17").Value2
$ws.Range("O16").Value = $ws.Range("O17").Value2
$ws.Range("P16").Value = 580
$ws.Range("Q16").Value = $ws.Range("Q17").Value2
$ws.Range("R16").Value = $ws.Range("R17").Value2
